$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Formatting applied to the numeric cells (B1, A2): bold font, thin box border,
# horizontally centered and top-aligned vertically.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108   # xlCenter
$b1.VerticalAlignment = -4160     # xlTop
$b1.Borders.LineStyle = 1         # xlContinuous
$b1.Borders.Weight = 2            # xlThin

# Re-use the exact same style for A2 (instead of re-deriving it step by step,
# which would otherwise create extra, unused intermediate cell styles).
$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
